$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The logging tool that produced this workbook stamps the current
# UserName/WorkGroup into row 2 (A2/C2) every time it runs, appending a
# fresh timestamped entry to the shared-string history.  This commit
# records a new run, so update A2 and C2 to the newly generated values.
$ws.Range("A2").Value = "UserName1540893075711"
$ws.Range("C2").Value = "WorkGroup1540893191201"
